# Update the "想去人数" (interested-people count) figures in column F
# for the two worksheets that carry this exhibition data: "展览" and
# "全部类型". Both sheets contain identical rows 2-9 for these events,
# so the same six cells are updated on each sheet.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2" = 340
    "F4" = 500
    "F5" = 4998
    "F7" = 627
    "F8" = 293
    "F9" = 755
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
